# Apply "Change data input. Interpolate between high and low energy demand
# value for average number of occupants" edit to the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the annual electrical energy demand values (column C).
# The specific demand (column D) is computed via formulas already present
# in the sheet (C/B), so it recalculates automatically.
$ws.Range("C2").Value = 2000
$ws.Range("C3").Value = 3000
$ws.Range("C4").Value = 3500
$ws.Range("C5").Value = 3800
$ws.Range("C6").Value = 4800
$ws.Range("C7").Value = 1200
$ws.Range("C8").Value = 2000
$ws.Range("C9").Value = 2400
$ws.Range("C10").Value = 2800
$ws.Range("C11").Value = 3500

# Add reference information below the table.
$ws.Range("A13").Value = "Ref:"
$ws.Range("B13").Value = "https://www.mieterbund.de/index.php?eID=tx_nawsecuredl&u=0&g=0&t=1496495412&hash=624f834b069c77c42b4a96024b56a6944f10bce4&file=fileadmin/pdf/Stromspiegel/Stromspiegel-2017_Tabellen.pdf"
$ws.Range("B14").Value = "Kategorie C"

$ws.Range("B14").Select()

$wb.Save()
